$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "SynchronizeLogNewEntries"

$ws.Range("A1").Value = "WorkbookName"
$ws.Range("A2").Value = "Sets-DemoModels.xlsx"

$ws.Range("B1").Value = "SheetName"
$ws.Range("B2").Value = "-"

$ws.Range("C1").Value = "Type"
$ws.Range("C2").Value = "SetRules"

$ws.Range("D1").Value = "Scenario"
$ws.Range("D2").Value = "s-DemoModels"

$ws.Range("E1").Value = "Attribute"
$ws.Range("E2").Value = "-"

$ws.Range("F1").Value = "MessageCategory"
$ws.Range("F2").Value = "warning"

$ws.Range("G1").Value = "Message"
$ws.Range("G2").Value = "The following sets did not generate any records: NRG_NUK ,NRG_PP ,NRG_SOLID ,PP_NUCLEAR"

$ws.Range("H1").Value = "CellAddress"
$ws.Range("H2").Value = "-"

$ws.Range("I1").Value = "TagType"
$ws.Range("I2").Value = "-"

$ws.Range("J1").Value = "ProcessFilter"
$ws.Range("J2").Value = "-"

$ws.Range("K1").Value = "CommodityFilter"
$ws.Range("K2").Value = "-"
